$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.249.42"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "3.202.27"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.13"
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.37"
$ws.Range("E6").Value = "  +1.13%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.201.01"
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.67"
$ws.Range("E11").Value = "  -3.50%  "
$ws.Range("E12").Value = "  -3.07%  "
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "38.42"
$ws.Range("E14").Value = "  -2.24%  "
$ws.Range("D15").Value = "3.728.49"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "66.382.19"
$ws.Range("E17").Value = "  -2.74%  "
$ws.Range("D18").Value = "3.208.27"
$ws.Range("E18").Value = "  +0.64%  "
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "506.93"
$ws.Range("E20").Value = "  -2.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.31"
$ws.Range("E21").Value = "  -0.93%  "
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.00"
$ws.Range("E23").Value = "  -1.60%  "
$ws.Range("E24").Value = "  -2.48%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.13"
$ws.Range("E25").Value = "  -1.20%  "
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.99"
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.05"
$ws.Range("E28").Value = "  -2.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.36"
$ws.Range("E29").Value = "  +0.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.130"
$ws.Range("E30").Value = "  +44.10%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.99"
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("B32").Value = "Stacks"
$ws.Range("C32").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.92"
$ws.Range("E32").Value = "  -0.80%  "
$ws.Range("E33").Value = "  -0.51%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  -4.65%  "
$ws.Range("E36").Value = "  -1.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.41"
$ws.Range("E37").Value = "  +0.75%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "500.54"
$ws.Range("E38").Value = "  -2.09%  "
$ws.Range("D39").Value = "0.0₃0770"
$ws.Range("E39").Value = "  +14.38%  "
$ws.Range("B40").Value = "Kaspa"
$ws.Range("C40").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.131"
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0420"
$ws.Range("E41").Value = "  -1.35%  "
$ws.Range("E42").Value = "  +4.88%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.71"
$ws.Range("E43").Value = "  -2.38%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.296"
$ws.Range("E44").Value = "  -1.96%  "
$ws.Range("D45").Value = "2.908.72"
$ws.Range("E45").Value = "  +0.22%  "
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("E47").Value = "  -1.41%  "
$ws.Range("E48").Value = "  +1.59%  "
$ws.Range("E50").Value = "  -0.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "122.65"
$ws.Range("E51").Value = "  +0.35%  "
